$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Nr. Topics" column (F3:F7) to H3:H7
$ws.Range("H3").Value = $ws.Range("F3").Value()
$ws.Range("H4").Value = $ws.Range("F4").Value()
$ws.Range("H5").Value = $ws.Range("F5").Value()
$ws.Range("H6").Value = $ws.Range("F6").Value()
$ws.Range("H7").Value = $ws.Range("F7").Value()

# Clear the old F column data (rows 4-7), and set new header for F3
$ws.Range("F4").Value = $null
$ws.Range("F5").Value = $null
$ws.Range("F6").Value = $null
$ws.Range("F7").Value = $null
$ws.Range("F3").Value = "ncomp"

# Update selection to F4
$ws.Range("F4").Select()
